$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had two scales ("Initiate Scale" at row 8 and "Organization Scale"
# at row 12) that were never filled in with any Pre Baseline Phase data.
# Those placeholder rows are removed, and a brand-new "Pre Experimental
# Phase" column of results (column C) is added for every remaining scale.
$ws.Rows("12").Delete()
$ws.Rows("8").Delete()

# Column A - scale/index labels in their final row order
$ws.Range("A2").Value = "Inhibit Scale"
$ws.Range("A3").Value = "Self-Monitor Scale"
$ws.Range("A4").Value = "Behavior Regulation Index"
$ws.Range("A5").Value = "Shift Scale"
$ws.Range("A6").Value = "Emotional Control Scale"
$ws.Range("A7").Value = "Emotional Regulation Index"
$ws.Range("A8").Value = "Task Completion Scale"
$ws.Range("A9").Value = "Working Memory Scale"
$ws.Range("A10").Value = "Plan/Organize Scale"
$ws.Range("A11").Value = "Cognitive Regulation Index"
$ws.Range("A12").Value = "Global Executive Composite "

# Column B - Pre Baseline Phase scores (unchanged values, re-asserted in
# their new row positions)
$ws.Range("B2").Value = 64
$ws.Range("B3").Value = 76
$ws.Range("B4").Value = 71
$ws.Range("B5").Value = 64
$ws.Range("B6").Value = 71
$ws.Range("B7").Value = 69
$ws.Range("B8").Value = 78
$ws.Range("B9").Value = 72
$ws.Range("B10").Value = 60
$ws.Range("B11").Value = 71
$ws.Range("B12").Value = 72

# Column C - new Pre Experimental Phase scores
$ws.Range("C2").Value = 79
$ws.Range("C3").Value = 72
$ws.Range("C4").Value = 79
$ws.Range("C5").Value = 71
$ws.Range("C6").Value = 64
$ws.Range("C7").Value = 69
$ws.Range("C8").Value = 82
$ws.Range("C9").Value = 84
$ws.Range("C10").Value = 73
$ws.Range("C11").Value = 82
$ws.Range("C12").Value = 80

# Match the author's final cursor position
$ws.Range("C13").Select()
